$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 38 (A38 = icsdadultosa_psg5events): set poordi3 mapping
$ws.Range("B38").Value = "poordi3"
$ws.Range("C38").Value = "Obstructive apnea (all desaturations) hypopnea (3% desaturation) index"
$ws.Range("D38").Value = "x"
$ws.Range("D38").ClearContents()
$ws.Range("D38").Style = "Normal"

# Row 39 (A39 = icsdadultosa_psg15events): set poordi3 mapping
$ws.Range("B39").Value = "poordi3"
$ws.Range("C39").Value = "Obstructive apnea (all desaturations) hypopnea (3% desaturation) index"
$ws.Range("D39").Value = "x"
$ws.Range("D39").ClearContents()
$ws.Range("D39").Style = "Normal"

# Update selection to match the newly-edited rows
$ws.Range("B39:D39").Select()
